$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update status of "Logging-System" (row 13) from "in Arbeit" to "done"
$ws.Range("B13").Value = "done"
$ws.Range("B2").Copy()
$ws.Range("B13").PasteSpecial(-4122)

# Add new TODO entry in row 16
$ws.Range("A16").Value = "Konstanten in JSON-File"
$ws.Range("B16").Value = "offen"
$ws.Range("B4").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("C16").Value = "Jonas"

# The "Neutral" (in Arbeit) cell style is no longer used anywhere in the
# workbook now that row 13 has been switched to "done" - remove it.
$wb.Styles.Item("Neutral").Delete()

# Leave the cursor on the cell that was last edited.
$ws.Range("B13").Select()
